$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.401.63"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.847.16"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9986"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.02"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6294"
$ws.Range("E6").Value = "  -3.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07606"
$ws.Range("E8").Value = "  +1.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2972"
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.46"
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("D11").Value = "2.200.76"
$ws.Range("E11").Value = "  +18.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07715"
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6877"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.986"
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").Value = "2.313.16"
$ws.Range("E15").Value = "  +9.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.90"
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009887"
$ws.Range("E17").Value = "  +4.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.167"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").Value = "29.411.96"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "231.77"
$ws.Range("E20").Value = "  -2.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.51"
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.608"
$ws.Range("E23").Value = "  -1.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9996"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.74"
$ws.Range("E25").Value = "  -1.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1391"
$ws.Range("E26").Value = "  -1.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.475"
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.69"
$ws.Range("E28").Value = "  -0.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.471"
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("E30").Value = "  -4.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.252"
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.122"
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.017"
$ws.Range("E33").Value = "  -1.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.867"
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.159"
$ws.Range("E35").Value = "  -2.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7172"
$ws.Range("E36").Value = "  -1.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.594"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "1.248.06"
$ws.Range("E38").Value = "  +3.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.791"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("E40").Value = "  +1.09%  "
$ws.Range("D41").Value = "2.267.80"
$ws.Range("E41").Value = "  +12.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9073"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.099"
$ws.Range("E43").Value = "  -2.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9993"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "67.43"
$ws.Range("E45").Value = "  +1.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.13"
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.299"
$ws.Range("E47").Value = "  -1.85%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000118"
$ws.Range("E48").Value = "  -4.79%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.174"
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4012"
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.699"
$ws.Range("E51").Value = "  +2.98%  "
